# Updates cryptos list: refresh Price (D) / Volume(1h) (E) figures,
# and shift Coin/Link/Price/Volume rows 20-51 up by one ranking slot
# (a new row, Elrond, enters at the bottom, row 51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "31.166.09"
$ws.Range("E2").Value = "  +2.23%  "

$ws.Range("D3").Value = "1.970.55"
$ws.Range("E3").Value = "  +3.21%  "

$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.30%  "

$ws.Range("D5").Value = "'248.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.04%  "

$ws.Range("D6").Value = "'1.004"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.40%  "

$ws.Range("D7").Value = "'0.4895"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.18%  "

$ws.Range("D8").Value = "'44.93"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.12%  "

$ws.Range("D9").Value = "'0.2953"
$ws.Range("D9").Style = "Normal"

$ws.Range("D10").Value = "'0.06834"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.43%  "

$ws.Range("D11").Value = "'19.19"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.68%  "

$ws.Range("D12").Value = "'107.44"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.55%  "

$ws.Range("D13").Value = "1.963.84"
$ws.Range("E13").Value = "  +2.78%  "

$ws.Range("D14").Value = "'0.07789"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.05%  "

$ws.Range("D15").Value = "'5.442"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.16%  "

$ws.Range("D16").Value = "'0.7072"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.05%  "

$ws.Range("D17").Value = "'285.60"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.31%  "

$ws.Range("D18").Value = "31.183.85"
$ws.Range("E18").Value = "  +2.23%  "

$ws.Range("D19").Value = "'13.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.82%  "

$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "'0.000007757"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.66%  "

$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "2.222.51"
$ws.Range("E21").Value = "  +2.72%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'5.626"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.63%  "

$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").Value = "'1.005"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.43%  "

$ws.Range("B24").Value = "BinanceUSD"
$ws.Range("C24").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D24").Value = "'1.003"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.29%  "

$ws.Range("B25").Value = "Chainlink"
$ws.Range("C25").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D25").Value = "'6.650"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.31%  "

$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "'10.02"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.32%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "'170.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.99%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'20.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.26%  "

$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").Value = "'2.192"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.17%  "

$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").Value = "'0.1067"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.28%  "

$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").Value = "'1.446"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.21%  "

$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'4.829"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +19.34%  "

$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "'4.524"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +9.74%  "

$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.05084"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.30%  "

$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'0.7711"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.11%  "

$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "'1.171"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.76%  "

$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "'2.740"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.92%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.02047"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.82%  "

$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "'2.732"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.93%  "

$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'6.465"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +12.10%  "

$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "'2.128"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.22%  "

$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "'74.09"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.47%  "

$ws.Range("D43").Value = "'0.8867"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.00%  "

$ws.Range("D44").Value = "'110.09"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.78%  "

$ws.Range("D45").Value = "'0.4475"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.44%  "

$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value = "'1.004"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.45%  "

$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").Value = "'7.531"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.14%  "

$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "'988.73"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +17.66%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'9.463"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.68%  "

$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "'0.1272"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.92%  "

$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").Value = "'36.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.71%  "
